$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'pants for men sport',
    'mens pouch leggings',
    'spandex leggings boys',
    'rodilleras de basketball',
    'knee sleeves wrestling',
    'knee pads workout',
    'baseball catcher leg guards adult',
    '5 inch seam shorts men',
    'youth softball pants for girls',
    'mens compression workout pants',
    'compression shorts for men',
    'work knee pads for men',
    'padded shorts snowboarding',
    'youth girls yoga pants',
    'youth hockey padded shorts',
    'knee protector sports',
    'kneeling pad construction',
    'knee pad for work',
    'basketballs under',
    'boys baseball pants size 6',
    'sliding short',
    'black knee pads volleyball girls',
    'mens 3/4 pants',
    'football leg sleeves for men',
    'black basketball shorts men',
    'compression football shorts',
    'running leggings men',
    'football leggings boys',
    'baseball leg guards',
    'protective basketball',
    'mountain bike knee pads',
    'weightlifting shorts men',
    'shorts for men basketball',
    'wrestling shorts for boys',
    'baseball items for men',
    'knee pads for biking',
    'eva foam knee pads',
    'mens compression running tights',
    'womens softball pants black',
    'waist guard',
    'bump pads',
    'mens fitness pants',
    'cycling pants for men padded',
    'knee pads work',
    'youth football girdle',
    'bjj knee sleeves',
    'volleyball kneepads black',
    'non slip knee pads',
    'short baseball',
    'knee sleeves for wrestling',
    'knee pads',
    'compression sleeve youth baseball',
    'tights for football',
    'soccer pad',
    'cycling pants for men',
    'dry fit leggings men',
    'compression calf sleeve men basketball',
    'long compression shorts men',
    'compression shorts long men',
    'compression pants and tops for men',
    'leggings knee length',
    'mens softball gear',
    'yoga after knee replacement',
    'wrestling knee sleeve youth',
    'tights compression',
    'mens compression pants pack',
    'boys running pants',
    'knee pads thigh support',
    'youth baseball pants black',
    'knee pads biking adult',
    'youth boys leggings',
    'adult pants',
    'youth baseball compression sleeves',
    'calf sleeves for men football',
    'padded knee sleeve',
    'knee pad exercise',
    'recovery pants men',
    'mens tight',
    'mens outdoor basketball',
    'soccer pants youth',
    'protective shorts',
    'baseball hand guard',
    'bee pants',
    'mens protective pads',
    'mens shorts above knee',
    'basketball knee support for men',
    'compression pads for surgery',
    'snowboarding pants boys',
    'basketball pants for girls',
    'youth baseball pants girls',
    'knee pads for exercise',
    'girls compression knee sleeves',
    'men yoga pant',
    'yoga pants mens black',
    'basketball aids',
    'knee compression sleeve volleyball',
    'impact advanced recovery',
    'best basketball',
    'baseball compression sleeve youth',
    'compression knee sleeve padded',
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
